$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Anselmo-Gestão integr"

$ws.Range("B3").Value = "Cleidson-Circuitos elétri"
$ws.Range("C3").Value = "Anselmo-Gestão integr"

$ws.Range("B4").Value = "Cleidson-Circuitos elétri"
$ws.Range("C4").Value = "[-, Joel L.-Tecnologia dos Materiais.]"
$ws.Range("E4").Value = "[-, -, -, Valmir-Metrologia]"

$ws.Range("B6").Value = "Cleidson-Circuitos elétri"
$ws.Range("C6").Value = "[-, Joel L.-Tecnologia dos Materiais.]"
$ws.Range("E6").Value = "[-, -, -, Valmir-Metrologia]"

$ws.Range("B7").Value = "Cleidson-Circuitos elétri"
$ws.Range("C7").Value = "André Guimarães-Desenho Técn"
$ws.Range("E7").Value = "[-, -, -, Valmir-Metrologia]"
$ws.Range("F7").Value = "-"

$ws.Range("E8").Value = "[-, -, -, Valmir-Metrologia]"
$ws.Range("F8").Value = "-"
